$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "Beteckning" (id) for every data row; use it to find
# the last populated row so the update covers the whole data table
# (rows 2..lastRow), mirroring the "Förändrad" (Changed) column C update.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$newDate = 45202

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
